# "Generate Report for Handoff"
#
# The dd442c70-...md file has finished translation handoff generation for
# both zh-cn and de-de. Update its status from "In Translation" to
# "Ready for handoff", bump its priority from "ht" (human translation) to
# "mt" (machine translation), and record the new handoff timestamps across
# the Overview rollup sheet and the two per-locale detail sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (row 3 = dd442c70-...md) ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-08-31 16:16:16"

# ---- zh-cn detail sheet (row 3 = dd442c70-...md) ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "mt"
$zh.Range("H3").Value = "2016-08-31 16:16:10"

# ---- de-de detail sheet (row 3 = dd442c70-...md) ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "mt"
$de.Range("H3").Value = "2016-08-31 16:16:16"

# ---- Widen the Status-ish columns that now hold the longer
#      "Ready for handoff" text (Overview E:F, zh-cn/de-de C) ----
$targetWidth = 16.333333333333332
$ovw.Columns.Item(5).ColumnWidth = $targetWidth
$ovw.Columns.Item(6).ColumnWidth = $targetWidth
$zh.Columns.Item(3).ColumnWidth = $targetWidth
$de.Columns.Item(3).ColumnWidth = $targetWidth
